# Apply text replacements per the diff: date update and equation value updates.
$d = $word.ActiveDocument

$replacements = @(
    @('2024-09-08 Sunday', '2024-09-09 Monday'),
    @('26+36=62', '85-16=69'),
    @('64+4=68', '61+19=80'),
    @('44+54=98', '41-20=21'),
    @('98-12=86', '93-62=31'),
    @('91-27=64', '53-17=36'),
    @('44+24=68', '81-55=26'),
    @('44-32=12', '43+8=51'),
    @('37+18=55', '91-59=32'),
    @('66+25=91', '17+79=96'),
    @('4+28=32', '59-23=36'),
    @('40+51=91', '73-48=25'),
    @('13+28=41', '59-43=16'),
    @('66+3=69', '76-7=69'),
    @('21+29=50', '11+30=41'),
    @('66-65=1', '65+18=83'),
    @('42-42=0', '7+32=39'),
    @('92-6=86', '67-19=48'),
    @('1+77=78', '86-74=12'),
    @('38+2=40', '51+42=93'),
    @('8+75=83', '61-44=17'),
    @('28-10=18', '81-72=9'),
    @('98-9=89', '32-16=16'),
    @('36+7=43', '45+12=57'),
    @('33-22=11', '90-29=61'),
    @('53+44=97', '44+48=92'),
    @('2+61=63', '93-43=50'),
    @('97-76=21', '10+29=39'),
    @('5+93=98', '79-4=75'),
    @('71-20=51', '12+44=56'),
    @('8+47=55', '43+4=47'),
    @('89-30=59', '27+4=31'),
    @('38+26=64', '9+71=80'),
    @('54+33=87', '38-34=4'),
    @('25-15=10', '0+35=35'),
    @('20+54=74', '22+0=22'),
    @('62-6=56', '90-1=89'),
    @('37-16=21', '0+57=57'),
    @('31-3=28', '9+81=90'),
    @('20+75=95', '25-8=17'),
    @('13+17=30', '68-53=15'),
    @('18+23=41', '4+0=4'),
    @('81-46=35', '33+42=75'),
    @('76-69=7', '65-61=4'),
    @('1+32=33', '41-10=31'),
    @('51+41=92', '54+12=66'),
    @('39-15=24', '1+56=57'),
    @('81-61=20', '73-28=45'),
    @('10+45=55', '96-38=58'),
    @('48+41=89', '6+55=61'),
    @('58+17=75', '42+39=81'),
    @('41+0=41', '56-46=10'),
    @('84-62=22', '76+17=93'),
    @('51-0=51', '5+46=51'),
    @('40-5=35', '5+76=81'),
    @('9+65=74', '85-63=22'),
    @('10+51=61', '90-85=5'),
    @('36+28=64', '24-15=9'),
    @('83+7=90', '52-44=8'),
    @('11+49=60', '94-74=20'),
    @('58+29=87', '56-35=21'),
    @('92-13=79', '16+23=39'),
    @('91-66=25', '50-6=44'),
    @('41-25=16', '31+67=98'),
    @('81+11=92', '35+18=53'),
    @('8+25=33', '43+11=54'),
    @('23+18=41', '91-68=23'),
    @('82-48=34', '81+14=95'),
    @('86-52=34', '95+1=96'),
    @('9+39=48', '87+11=98'),
    @('94-15=79', '75-32=43'),
    @('90+0=90', '89-19=70'),
    @('95-25=70', '74+21=95'),
    @('63-44=19', '19+8=27'),
    @('72+12=84', '65-35=30'),
    @('46+25=71', '17+2=19'),
    @('88-48=40', '78-15=63'),
    @('54-44=10', '63-45=18'),
    @('85-41=44', '60-48=12'),
    @('83+2=85', '15+11=26'),
    @('22-11=11', '89+8=97'),
    @('51-41=10', '52-50=2'),
    @('42+45=87', '26+29=55'),
    @('18+64=82', '20+65=85'),
    @('93-50=43', '44+52=96'),
    @('55+40=95', '11-1=10'),
    @('57-52=5', '24+12=36'),
    @('66-33=33', '72+13=85'),
    @('67-13=54', '40-24=16'),
    @('60+38=98', '7+75=82'),
    @('17+18=35', '66+10=76'),
    @('99-97=2', '78-50=28'),
    @('62-51=11', '6+93=99'),
    @('4+49=53', '21+61=82'),
    @('26+8=34', '39+31=70'),
    @('0+62=62', '44-41=3'),
    @('80-72=8', '24+36=60'),
    @('15-10=5', '80-41=39'),
    @('67+13=80', '59-45=14'),
    @('80-58=22', '23+12=35'),
    @('87+0=87', '29+26=55')
)

$count = 0
foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if ($found) {
        $count = $count + 1
    } else {
        Write-Output "NOT FOUND: $oldText"
    }
}

Write-Output "Replaced $count of $($replacements.Count)"
